# Swap the data (all columns except the running index in column A) between
# pairs of rows in the "South Korea K3 League" sheet, as described by the
# commit's diff (rows were reordered / re-paired).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose contents (columns B..AD) must be exchanged.
$rowPairs = @(
    @(14, 15),
    @(35, 36),
    @(40, 41),
    @(49, 50),
    @(93, 94),
    @(118, 119),
    @(124, 125),
    @(126, 127),
    @(200, 201)
)

$firstCol = 2   # column B
$lastCol  = 30  # column AD

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Capture every value of both rows first so reading is not affected by
    # the writes we are about to perform.
    $row1Values = @{}
    $row2Values = @{}

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $row1Values[$col] = $ws.Cells.Item($r1, $col).Value2
        $row2Values[$col] = $ws.Cells.Item($r2, $col).Value2
    }

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $newVal1 = $row2Values[$col]
        $newVal2 = $row1Values[$col]

        if ($newVal1 -eq $null) {
            $cell1.ClearContents()
        } else {
            $cell1.Value = $newVal1
        }

        if ($newVal2 -eq $null) {
            $cell2.ClearContents()
        } else {
            $cell2.Value = $newVal2
        }
    }
}
